$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook's data table (rows 16-34) lists the "Valor Mora" (overdue
# amount) per worker per "Periodo Mora" (arrears period, format YYMM).
# This update re-sorts the table by Periodo Mora ascending (it previously
# ran from the most recent period down to the oldest) and folds in the
# first batch of new account-statement rows for two additional workers
# (NATALIA ANDREA RICARDO MENA / CC 1007980980 and JOSE DANIEL VILLA
# VILLARREAL / CC 1007229509), each owing against period 2202 / 2204
# respectively - the same periods RICARDO ANDRES BARRIOS MONTES already
# has rows for, so they land right after his matching row once sorted.
#
# Only the cell VALUES move; the per-row formatting (s="..") already in
# place on the sheet is left exactly as-is, so we write straight into the
# fixed row positions instead of doing a Range.Sort (which would also drag
# the formatting along with the data).

$rows = @(
    @("CC", "1143384213", "RICARDO ANDRES BARRIOS MONTES", "2202", 80000),
    @("CC", "1007980980", "NATALIA ANDREA RICARDO MENA",   "2202", 80000),
    @("CC", "1143384213", "RICARDO ANDRES BARRIOS MONTES", "2203", 80000),
    @("CC", "1143384213", "RICARDO ANDRES BARRIOS MONTES", "2204", 80000),
    @("CC", "1007229509", "JOSE DANIEL VILLA VILLARREAL",  "2204", 53334),
    @("CC", "1143384213", "RICARDO ANDRES BARRIOS MONTES", "2205", 80000),
    @("CC", "1143384213", "RICARDO ANDRES BARRIOS MONTES", "2206", 80000),
    @("CC", "1143384213", "RICARDO ANDRES BARRIOS MONTES", "2207", 80000),
    @("CC", "1143384213", "RICARDO ANDRES BARRIOS MONTES", "2208", 80000),
    @("CC", "1143384213", "RICARDO ANDRES BARRIOS MONTES", "2209", 80000),
    @("CC", "1143384213", "RICARDO ANDRES BARRIOS MONTES", "2210", 80000),
    @("CC", "1143384213", "RICARDO ANDRES BARRIOS MONTES", "2211", 80000),
    @("CC", "1143384213", "RICARDO ANDRES BARRIOS MONTES", "2212", 80000),
    @("CC", "1143384213", "RICARDO ANDRES BARRIOS MONTES", "2301", 80000),
    @("CC", "1143384213", "RICARDO ANDRES BARRIOS MONTES", "2302", 80000),
    @("CC", "1143384213", "RICARDO ANDRES BARRIOS MONTES", "2303", 80000),
    @("CC", "1143384213", "RICARDO ANDRES BARRIOS MONTES", "2304", 80000),
    @("CC", "1143384213", "RICARDO ANDRES BARRIOS MONTES", "2305", 80000),
    @("CC", "1143384213", "RICARDO ANDRES BARRIOS MONTES", "2306", 58667)
)

$r = 16
foreach ($row in $rows) {
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $r = $r + 1
}
